# "se crea el programa copiar_info_reunion"
# Rename the single worksheet and move the active selection, matching the
# saved state captured in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet1 -> Info-reunión
$ws.Name = "Info-reunión"

# Move / persist the active selection to B19
$ws.Activate()
$ws.Range("B19").Select()
